# Auto-generated edit script applying "Add data for 2023-11-03" changes
# to output/violent-crime-full-year.xlsx (violent crime dataset workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6473
$ws.Range("J3").Value = 6875
$ws.Range("H4").Value = 1705
$ws.Range("J4").Value = 1488
$ws.Range("J5").Value = 530
$ws.Range("J6").Value = 9098
$ws.Range("H7").Value = 26016
$ws.Range("J7").Value = 24464

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("J2").Value = 20
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 418
$ws.Range("J7").Value = 1542

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J6").Value = 131
$ws.Range("J7").Value = 484

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 261
$ws.Range("J3").Value = 363
$ws.Range("J7").Value = 1098

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 121
$ws.Range("J7").Value = 354

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 255
$ws.Range("J7").Value = 750

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 178
$ws.Range("J4").Value = 26
$ws.Range("J6").Value = 218
$ws.Range("J7").Value = 608

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 104
$ws.Range("J5").Value = 11
$ws.Range("J7").Value = 375

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 112
$ws.Range("J7").Value = 712
$ws.Range("J8").Value = 1542
$ws.Range("J9").Value = 129
$ws.Range("J11").Value = 413
$ws.Range("J12").Value = 50
$ws.Range("J15").Value = 284
$ws.Range("J18").Value = 209
$ws.Range("J19").Value = 718
$ws.Range("J20").Value = 507
$ws.Range("J23").Value = 225
$ws.Range("J24").Value = 77
$ws.Range("J29").Value = 1333
$ws.Range("J31").Value = 233
$ws.Range("J32").Value = 40
$ws.Range("J33").Value = 1098
$ws.Range("J36").Value = 333
$ws.Range("J37").Value = 750
$ws.Range("J42").Value = 1051
$ws.Range("J47").Value = 184
$ws.Range("J51").Value = 301
$ws.Range("J52").Value = 618
$ws.Range("J54").Value = 466
$ws.Range("J55").Value = 372
$ws.Range("J56").Value = 34
$ws.Range("J57").Value = 109
$ws.Range("I63").Value = 250
$ws.Range("J64").Value = 162
$ws.Range("J65").Value = 608
$ws.Range("J67").Value = 926
$ws.Range("J68").Value = 53
$ws.Range("J69").Value = 53
$ws.Range("J71").Value = 78
$ws.Range("J73").Value = 235
$ws.Range("J76").Value = 365
$ws.Range("J78").Value = 291
$ws.Range("H79").Value = 768
$ws.Range("I79").Value = 748
$ws.Range("J79").Value = 688
$ws.Range("J83").Value = 484
$ws.Range("J85").Value = 1020
$ws.Range("J86").Value = 158
$ws.Range("J88").Value = 250
$ws.Range("J89").Value = 318
$ws.Range("J94").Value = 257
$ws.Range("J95").Value = 354
$ws.Range("J96").Value = 268
$ws.Range("J99").Value = 375
$ws.Range("H101").Value = 26016
$ws.Range("J101").Value = 24464

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 85
$ws.Range("J7").Value = 233

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 234
$ws.Range("J3").Value = 343
$ws.Range("J6").Value = 258
$ws.Range("J7").Value = 926

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 113
$ws.Range("J6").Value = 219
$ws.Range("J7").Value = 466

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 406
$ws.Range("J4").Value = 70
$ws.Range("J7").Value = 1333

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 168
$ws.Range("J6").Value = 280
$ws.Range("J7").Value = 718

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 61
$ws.Range("J3").Value = 76
$ws.Range("J7").Value = 365

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 204
$ws.Range("J6").Value = 560
$ws.Range("J7").Value = 1051

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 291

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J6").Value = 207
$ws.Range("J7").Value = 372

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 225

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 40
$ws.Range("J4").Value = 43
$ws.Range("H7").Value = 768
$ws.Range("I7").Value = 748
$ws.Range("J7").Value = 688

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J2").Value = 45
$ws.Range("J6").Value = 55
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 174
$ws.Range("J7").Value = 507

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 56
$ws.Range("J7").Value = 209

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 106
$ws.Range("J7").Value = 333

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J6").Value = 231
$ws.Range("J7").Value = 712

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 142
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 87
$ws.Range("J7").Value = 184

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 65
$ws.Range("J6").Value = 121
$ws.Range("J7").Value = 284

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J6").Value = 185
$ws.Range("J7").Value = 413

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J3").Value = 41
$ws.Range("J7").Value = 129

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J3").Value = 61
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 235

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J6").Value = 120
$ws.Range("J7").Value = 250

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 92
$ws.Range("J7").Value = 318

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 86
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J6").Value = 120
$ws.Range("J7").Value = 301

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 360
$ws.Range("J6").Value = 295
$ws.Range("J7").Value = 1020

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 181
$ws.Range("J4").Value = 22
$ws.Range("J6").Value = 263
$ws.Range("J7").Value = 618

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 112

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 50
